# WebForm User Assignment execution
# Update the phone-number (PN_Value) column F for rows 2-10 and the
# Match2UserPos value in AN2 to reflect the latest webform run data.
#
# NumberFormat "@" forces the numeric-looking digit strings to be stored
# as text (matching the original shared-string cell type) instead of
# being coerced to a number; ClearFormats() afterwards drops the
# temporary text format so the cell keeps the sheet's default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.ClearFormats()
}

Set-TextValue "F2"  "9840006389"
Set-TextValue "F3"  "9840024574"
Set-TextValue "F4"  "9840070213"
Set-TextValue "F5"  "9840080807"
Set-TextValue "F6"  "9840092307"
Set-TextValue "F7"  "9840085281"
Set-TextValue "F8"  "9840054735"
Set-TextValue "F9"  "9840059770"
Set-TextValue "F10" "9840070370"

Set-TextValue "AN2" "2"
